$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Testdata_path" column, making room
# for the new "Testdata_name" column. This preserves the existing B column
# (now shifted to C) in place, including its shared-string usage.
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Testdata_name"

# Fill down the "test" environment value for the newly added rows
$ws.Range("A3").Value = "test"
$ws.Range("A4").Value = "test"
$ws.Range("A5").Value = "test"
$ws.Range("A6").Value = "test"
$ws.Range("A7").Value = "test"

# Rows for NonOnco ManagePopulations data (LIVEHTA-1328)
$ws.Range("B2").Value = "nononcology_managepopulationdata"
$ws.Range("B3").Value = "managepopulation_additional_col_check"
$ws.Range("C2").Value = "\Testdata\Non_Oncology\DataFiles\ManagePopulations\NonOncologyManagePopulationsPage_Data.xlsx"
$ws.Range("C3").Value = "\Testdata\Non_Oncology\DataFiles\ManagePopulations\ManagePopulations_AdditionalCol_Check_Data.xlsx"

# Rows for Edit Population invalid-data scenarios
$ws.Range("B4").Value = "edit_ep_categorical_invaliddata"
$ws.Range("B5").Value = "edit_ep_continuous_invaliddata"
$ws.Range("B6").Value = "edit_ep_timetoevent_invaliddata"
$ws.Range("C4").Value = "\Testdata\Non_Oncology\DataFiles\ManagePopulations\EP_Categorical\EditPopulations_with_Invalid_Categorical_Data.xlsx"
$ws.Range("C5").Value = "\Testdata\Non_Oncology\DataFiles\ManagePopulations\EP_Continuous\EditPopulations_with_Invalid_Continuous_Data.xlsx"
$ws.Range("C6").Value = "\Testdata\Non_Oncology\DataFiles\ManagePopulations\EP_TimetoEvent\EditPopulations_with_Invalid_TimetoEvent_Data.xlsx"

# Row for NonOnco Import tool (LIVEHTA-1449)
$ws.Range("B7").Value = "nononcology_importtool"
$ws.Range("C7").Value = "\Testdata\Non_Oncology\DataFiles\ImportPublications\ImportPublicationsPage_Data.xlsx"

$ws.Range("A1:C7").EntireColumn.AutoFit()

$ws.Range("A7").Select()
